$wb = $excel.ActiveWorkbook

# --- Sheet 1: "展览" (first worksheet) ---
$ws1 = $wb.Worksheets.Item(1)

# Update the "想去人数" (want-to-go count) figures for existing rows.
$ws1.Range("F3").Value = 180
$ws1.Range("F4").Value = 773

# Append a new row (row 6) with a new event.
# Copy the bold/bordered/centered style used by the index column (A) from the
# row above, then set the new index value.
$ws1.Cells.Item(5,1).Copy()
$ws1.Cells.Item(6,1).PasteSpecial(-4122)
$ws1.Cells.Item(6,1).Value = 5

$c = $ws1.Cells.Item(6,2)
$c.NumberFormat = "@"
$c.Value = "2024-11-23"
$c.Style = "Normal"

$ws1.Cells.Item(6,3).Value = "南宁·黑塔利亚同人ONLY"
$ws1.Cells.Item(6,4).Value = "北湖北路48-5号(近北湖小区) 金御华尊国际大酒店"
$ws1.Cells.Item(6,5).Value = "2024.11.23 09:30-11.23 16:00"
$ws1.Cells.Item(6,6).Value = 1
$ws1.Cells.Item(6,7).Value = 58
$ws1.Cells.Item(6,8).Value = "https://show.bilibili.com/platform/detail.html?id=93131"
$ws1.Cells.Item(6,9).Value = "//i1.hdslb.com/bfs/openplatform/202409/22sR9Hkq1726988625762.jpeg"

# --- Sheet 4: "全部类型" (fourth worksheet) ---
$ws4 = $wb.Worksheets.Item(4)

# Update the "想去人数" figures for existing rows.
$ws4.Range("F4").Value = 180
$ws4.Range("F5").Value = 773

# Append a new row (row 7) with the same new event.
$ws4.Cells.Item(6,1).Copy()
$ws4.Cells.Item(7,1).PasteSpecial(-4122)
$ws4.Cells.Item(7,1).Value = 6

$c2 = $ws4.Cells.Item(7,2)
$c2.NumberFormat = "@"
$c2.Value = "2024-11-23"
$c2.Style = "Normal"

$ws4.Cells.Item(7,3).Value = "南宁·黑塔利亚同人ONLY"
$ws4.Cells.Item(7,4).Value = "北湖北路48-5号(近北湖小区) 金御华尊国际大酒店"
$ws4.Cells.Item(7,5).Value = "2024.11.23 09:30-11.23 16:00"
$ws4.Cells.Item(7,6).Value = 1
$ws4.Cells.Item(7,7).Value = 58
$ws4.Cells.Item(7,8).Value = "https://show.bilibili.com/platform/detail.html?id=93131"
$ws4.Cells.Item(7,9).Value = "//i1.hdslb.com/bfs/openplatform/202409/22sR9Hkq1726988625762.jpeg"
